$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.512.88'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '2.530.67'
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.01'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.12'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.588'
$ws.Range('E7').Value = '  +1.85%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.540'
$ws.Range('E9').Value = '  -2.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.65'
$ws.Range('E10').Value = '  -0.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('D14').Value = '2.920.22'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').Value = '2.527.29'
$ws.Range('E15').Value = '  -2.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.22'
$ws.Range('E16').Value = '  +5.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.862'
$ws.Range('E17').Value = '  -2.92%  '
$ws.Range('D18').Value = '42.537.43'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.94'
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.45'
$ws.Range('E21').Value = '  -2.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.16'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '251.68'
$ws.Range('E23').Value = '  -1.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.91'
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.02'
$ws.Range('E25').Value = '  -4.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.91'
$ws.Range('E26').Value = '  -6.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.33'
$ws.Range('E28').Value = '  +9.66%  '
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.17'
$ws.Range('E30').Value = '  +0.85%  '
$ws.Range('E31').Value = '  -1.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '155.17'
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.33'
$ws.Range('E33').Value = '  -1.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0791'
$ws.Range('E34').Value = '  -2.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.07'
$ws.Range('E35').Value = '  -5.31%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.63'
$ws.Range('E36').Value = '  -4.66%  '
$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.56'
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('E38').Value = '  +1.92%  '
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.95'
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.40'
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.86'
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.02'
$ws.Range('E43').Value = '  -2.42%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0301'
$ws.Range('E45').Value = '  -3.17%  '
$ws.Range('D46').Value = '2.045.47'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '84.70'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.94'
$ws.Range('E48').Value = '  -3.42%  '
$ws.Range('D49').Value = '2.779.65'
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '101.88'
$ws.Range('E51').Value = '  -4.13%  '
